# Test Data Added for Slovakia market
#
# Duplicate the "Portugal" sheet (same layout/column widths/styles) to
# create a new "Slovakia" sheet, update its market/user-story cells, clean
# up the row heights that came along with the copy, and leave the
# selection/active-sheet state the way Excel would after this edit.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Copy Portugal to a new sheet placed right after it.
$portugal.Copy($null, $portugal)

$sheets = $wb.Worksheets
$slovakia = $sheets.Item($sheets.Count)
$slovakia.Name = "Slovakia"

# Market name + user story / ticket reference for the new market.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3236"

# The copied rows 3:5 inherited Portugal's explicit row height; re-fit them
# back down to the sheet default.
$slovakia.Rows("3:5").AutoFit()

# Leave the new sheet active, with C11 selected.
$slovakia.Activate()
$slovakia.Range("C11").Select() | Out-Null

# Portugal is no longer the active tab; its old B2 selection becomes a
# full-column (Select All) selection instead.
$portugal.Cells.Select() | Out-Null
$slovakia.Activate()
